$wb = $excel.ActiveWorkbook

# --- Cart_Page sheet: update locator rows 3-6, remove the now-unused blank
#     rows 7-8 gap (shifting the filler block up to close it), and make this
#     sheet the active/selected tab ---
$cart = $wb.Worksheets.Item("Cart_Page")

$cart.Range("A3").Value = "cart preview dropdown"
$cart.Range("B3").Value = "class name"
$cart.Range("C3").Value = "previewCart"
$cart.Range("E3").Value = "to find frame"

$cart.Range("A4").Value = "drop down products"
$cart.Range("B4").Value = "class name"
$cart.Range("C4").Value = "previewCartItem"

$cart.Range("A5").Value = "update link elements"
$cart.Range("B5").Value = "class name"
$cart.Range("C5").Value = "miniCart-qty-update"
$cart.Range("E5").ClearContents()

$cart.Range("A6").Value = "drop down cart"
$cart.Range("B6").Value = "id"
$cart.Range("C6").Value = "cart-preview-dropdown"

# Close the empty gap between row 6 and the filler block (old rows 9-23),
# which pulls that filler block up to rows 7-21 and drops the trailing two
# rows (22-23) that fall off the end.
$cart.Range("A7:A8").EntireRow.Delete()

$cart.Activate()
$cart.Range("A6").Select()

# --- Menu sheet: no longer the active tab; selection becomes the whole row 5 ---
$menu = $wb.Worksheets.Item("Menu")
$menu.Rows.Item(5).Select()

$cart.Activate()
